$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "1.000", "245.50") must be
# forced to Text format first, otherwise Excel auto-converts them to numeric values
# (stripping the exact decimal/trailing-zero formatting, e.g. "1.000" -> 1).
$ws.Range("D2").Value = '26.533.70'
$ws.Range("E2").Value = '  +4.02%  '
$ws.Range("D3").Value = '1.740.35'
$ws.Range("E3").Value = '  +4.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.50'
$ws.Range("E5").Value = '  +5.05%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4798'
$ws.Range("E7").Value = '  +4.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2685'
$ws.Range("E8").Value = '  +4.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06247'
$ws.Range("E9").Value = '  +2.09%  '
$ws.Range("D10").Value = '1.741.01'
$ws.Range("E10").Value = '  +4.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07145'
$ws.Range("E11").Value = '  +2.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.77'
$ws.Range("E12").Value = '  +8.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6181'
$ws.Range("E13").Value = '  +9.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.531'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.23'
$ws.Range("E15").Value = '  +3.46%  '
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '26.537.02'
$ws.Range("E17").Value = '  +4.06%  '
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006897'
$ws.Range("E19").Value = '  +3.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.75'
$ws.Range("E20").Value = '  +3.81%  '
$ws.Range("D21").Value = '1.965.73'
$ws.Range("E21").Value = '  +4.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.594'
$ws.Range("E22").Value = '  +4.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.864'
$ws.Range("E23").Value = '  +1.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.354'
$ws.Range("E24").Value = '  +3.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.82'
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.38'
$ws.Range("E26").Value = '  +3.61%  '
$ws.Range("E27").Value = '  +6.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.423'
$ws.Range("E28").Value = '  +4.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '107.13'
$ws.Range("E29").Value = '  +3.02%  '
$ws.Range("E30").Value = '  +1.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.734'
$ws.Range("E31").Value = '  +3.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07862'
$ws.Range("E32").Value = '  +1.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04556'
$ws.Range("E33").Value = '  +6.70%  '
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9992'
$ws.Range("E35").Value = '  +6.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6353'
$ws.Range("E36").Value = '  +6.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9312'
$ws.Range("E37").Value = '  +1.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '112.98'
$ws.Range("E38").Value = '  +11.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.431'
$ws.Range("E39").Value = '  -1.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.972'
$ws.Range("E40").Value = '  +9.19%  '
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.751'
$ws.Range("E42").Value = '  +17.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01508'
$ws.Range("E43").Value = '  +3.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3903'
$ws.Range("E44").Value = '  +5.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.763'
$ws.Range("E45").Value = '  +10.80%  '
$ws.Range("E46").Value = '  +9.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05330'
$ws.Range("E47").Value = '  +1.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.880'
$ws.Range("E48").Value = '  +7.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.74'
$ws.Range("E49").Value = '  +3.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.259'
$ws.Range("E50").Value = '  +5.62%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3438'
$ws.Range("E51").Value = '  +4.79%  '
